$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 40 entirely (shifting rows 41-44 up), which removes the
# "Flood Monitoring | National Framework for NMS/NHS Services | <flood monitoring paragraph>" row.
$ws.Rows.Item(40).Delete()

# Update the active selection cursor similar to the target state.
$ws.Range("D2").Select()
